$wb = $excel.ActiveWorkbook

# --- Training ---
$ws = $wb.Worksheets.Item("Training")

$ws.Range("F2").Value = 456214
$ws.Range("G2").Value = 9.281697591145834
$ws.Range("H2").Value = 0
$ws.Range("F3").Value = 394400
$ws.Range("G3").Value = 8.024088541666666
$ws.Range("H3").Value = 0

# TOTAL row 4: replace formulas with plain values
$ws.Range("D4").Value = 786432
$ws.Range("E4").Value = 2359296
$ws.Range("F4").Value = 850614
$ws.Range("G4").Value = 8.65289306640625
$ws.Range("H4").Value = 0

# --- Testing ---
$ws = $wb.Worksheets.Item("Testing")

$ws.Range("F2").Value = 500943
$ws.Range("G2").Value = 10.19171142578125
$ws.Range("H2").Value = 0
$ws.Range("F3").Value = 478046
$ws.Range("G3").Value = 9.725870768229166
$ws.Range("H3").Value = 0
$ws.Range("F4").Value = 391854
$ws.Range("G4").Value = 7.9722900390625
$ws.Range("H4").Value = 0
$ws.Range("F5").Value = 455834
$ws.Range("G5").Value = 9.273966471354166
$ws.Range("H5").Value = 0
$ws.Range("F6").Value = 427699
$ws.Range("G6").Value = 8.701558430989584
$ws.Range("H6").Value = 0
$ws.Range("F7").Value = 341927
$ws.Range("G7").Value = 6.956522623697917
$ws.Range("H7").Value = 0
$ws.Range("F8").Value = 522907
$ws.Range("G8").Value = 10.63857014973958
$ws.Range("H8").Value = 0
$ws.Range("F9").Value = 456200
$ws.Range("G9").Value = 9.281412760416666
$ws.Range("H9").Value = 0
$ws.Range("F10").Value = 401282
$ws.Range("G10").Value = 8.164103190104166
$ws.Range("H10").Value = 0
$ws.Range("F11").Value = 483296
$ws.Range("G11").Value = 9.832682291666666
$ws.Range("H11").Value = 0
$ws.Range("F12").Value = 433346
$ws.Range("G12").Value = 8.816446940104166
$ws.Range("H12").Value = 0
$ws.Range("F13").Value = 376433
$ws.Range("G13").Value = 7.658548990885417
$ws.Range("H13").Value = 0
$ws.Range("F14").Value = 408138
$ws.Range("G14").Value = 8.3035888671875
$ws.Range("H14").Value = 0
$ws.Range("F15").Value = 472926
$ws.Range("G15").Value = 9.6217041015625
$ws.Range("H15").Value = 0
$ws.Range("F16").Value = 348655
$ws.Range("G16").Value = 7.093404134114583
$ws.Range("H16").Value = 0
$ws.Range("F17").Value = 404093
$ws.Range("G17").Value = 8.221293131510416
$ws.Range("H17").Value = 0
$ws.Range("F18").Value = 376638
$ws.Range("G18").Value = 7.6627197265625
$ws.Range("H18").Value = 0
$ws.Range("F19").Value = 560798
$ws.Range("G19").Value = 11.40946451822917
$ws.Range("H19").Value = 0
$ws.Range("F20").Value = 430179
$ws.Range("G20").Value = 8.75201416015625
$ws.Range("H20").Value = 0
$ws.Range("F21").Value = 524120
$ws.Range("G21").Value = 10.66324869791667
$ws.Range("H21").Value = 0
$ws.Range("F22").Value = 399513
$ws.Range("G22").Value = 8.12811279296875
$ws.Range("H22").Value = 0
$ws.Range("F23").Value = 440634
$ws.Range("G23").Value = 8.9647216796875
$ws.Range("H23").Value = 0

# TOTAL row 24: replace formulas with plain values
$ws.Range("D24").Value = 8650752
$ws.Range("E24").Value = 25952256
$ws.Range("F24").Value = 9635461
$ws.Range("G24").Value = 8.910634358723957
$ws.Range("H24").Value = 0
$ws.Range("I24:N24").ClearContents()

# --- All Images ---
$ws = $wb.Worksheets.Item("All Images")

$ws.Range("F2").Value = 456214
$ws.Range("G2").Value = 9.281697591145834
$ws.Range("H2").Value = 0
$ws.Range("F3").Value = 394400
$ws.Range("G3").Value = 8.024088541666666
$ws.Range("H3").Value = 0
$ws.Range("F4").Value = 500943
$ws.Range("G4").Value = 10.19171142578125
$ws.Range("H4").Value = 0
$ws.Range("F5").Value = 478046
$ws.Range("G5").Value = 9.725870768229166
$ws.Range("H5").Value = 0
$ws.Range("F6").Value = 391854
$ws.Range("G6").Value = 7.9722900390625
$ws.Range("H6").Value = 0
$ws.Range("F7").Value = 455834
$ws.Range("G7").Value = 9.273966471354166
$ws.Range("H7").Value = 0
$ws.Range("F8").Value = 427699
$ws.Range("G8").Value = 8.701558430989584
$ws.Range("H8").Value = 0
$ws.Range("F9").Value = 341927
$ws.Range("G9").Value = 6.956522623697917
$ws.Range("H9").Value = 0
$ws.Range("F10").Value = 522907
$ws.Range("G10").Value = 10.63857014973958
$ws.Range("H10").Value = 0
$ws.Range("F11").Value = 456200
$ws.Range("G11").Value = 9.281412760416666
$ws.Range("H11").Value = 0
$ws.Range("F12").Value = 401282
$ws.Range("G12").Value = 8.164103190104166
$ws.Range("H12").Value = 0
$ws.Range("F13").Value = 483296
$ws.Range("G13").Value = 9.832682291666666
$ws.Range("H13").Value = 0
$ws.Range("F14").Value = 433346
$ws.Range("G14").Value = 8.816446940104166
$ws.Range("H14").Value = 0
$ws.Range("F15").Value = 376433
$ws.Range("G15").Value = 7.658548990885417
$ws.Range("H15").Value = 0
$ws.Range("F16").Value = 408138
$ws.Range("G16").Value = 8.3035888671875
$ws.Range("H16").Value = 0
$ws.Range("F17").Value = 472926
$ws.Range("G17").Value = 9.6217041015625
$ws.Range("H17").Value = 0
$ws.Range("F18").Value = 348655
$ws.Range("G18").Value = 7.093404134114583
$ws.Range("H18").Value = 0
$ws.Range("F19").Value = 404093
$ws.Range("G19").Value = 8.221293131510416
$ws.Range("H19").Value = 0
$ws.Range("F20").Value = 376638
$ws.Range("G20").Value = 7.6627197265625
$ws.Range("H20").Value = 0
$ws.Range("F21").Value = 560798
$ws.Range("G21").Value = 11.40946451822917
$ws.Range("H21").Value = 0
$ws.Range("F22").Value = 430179
$ws.Range("G22").Value = 8.75201416015625
$ws.Range("H22").Value = 0
$ws.Range("F23").Value = 524120
$ws.Range("G23").Value = 10.66324869791667
$ws.Range("H23").Value = 0
$ws.Range("F24").Value = 399513
$ws.Range("G24").Value = 8.12811279296875
$ws.Range("H24").Value = 0
$ws.Range("F25").Value = 440634
$ws.Range("G25").Value = 8.9647216796875
$ws.Range("H25").Value = 0

# TOTAL row 26: replace formulas with plain values
$ws.Range("D26").Value = 9437184
$ws.Range("E26").Value = 28311552
$ws.Range("F26").Value = 10486075
$ws.Range("G26").Value = 8.889155917697481
$ws.Range("H26").Value = 0
$ws.Range("I26:N26").ClearContents()
